$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.342.71"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.880.03"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.84%  "
$ws.Range("E6").Value = "  -2.87%  "
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.64"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.357"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.26"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0743"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("E13").Value = "  +3.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.153.33"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.765"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.93"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.865.45"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.331.90"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.59"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0823"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "244.43"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.85"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("E24").Value = "  +11.12%  "
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.16"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.93"
$ws.Range("D27").ClearFormats()
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.29"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.30"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0593"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("E33").Value = "  -2.32%  "
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("E35").Value = "  -11.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.42"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -12.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.857"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("E38").Value = "  -4.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0727"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.38"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0219"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.86"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.306.96"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0799"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.36%  "
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.83"
$ws.Range("D49").ClearFormats()
$ws.Range("E50").Value = "  -4.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.056.43"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.33%  "
